$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIVB")

# Row 12 - Gross Margin
$ws.Range("D12").Value = 1.0106
$ws.Range("E12").Value = 0.9812
$ws.Range("F12").Value = 0.9606
$ws.Range("G12").Value = 0.9397

# Row 14 - EBT margin
$ws.Range("D14").Value = 0.4242
$ws.Range("E14").Value = 0.3712
$ws.Range("F14").Value = 0.3997
$ws.Range("G14").Value = 0.4561

# Row 15 - Net Profit Margin
$ws.Range("D15").Value = 0.2967
$ws.Range("E15").Value = 0.2595
$ws.Range("F15").Value = 0.281
$ws.Range("G15").Value = 0.322

# Row 16 - Free Cash Flow Margin
$ws.Range("D16").Value = 0.3925
$ws.Range("E16").Value = 0.3963
$ws.Range("F16").Value = 0.2598
$ws.Range("G16").Value = 0.3112

# Row 23 - Operating Cash Flow Margin
$ws.Range("D23").Value = 0.4197
$ws.Range("E23").Value = 0.4249
$ws.Range("F23").Value = 0.2842
$ws.Range("G23").Value = 0.3297
